# Regenerate the "K" column (G) values in the save_data sheet.
# These are simulation-derived strikeout ("K") counts that replace the
# previous "Strike#" derived values; std/mean and s_vals were recalculated
# upstream and the sheet here only stores the resulting K counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 2
    4  = 4
    5  = 11
    6  = 6
    7  = 2
    8  = 1
    9  = 9
    10 = 6
    11 = 7
    12 = 1
    13 = 6
    14 = 5
    15 = 8
    16 = 4
    17 = 2
    18 = 7
    19 = 7
    20 = 10
    21 = 5
    22 = 7
    23 = 4
    24 = 2
    25 = 5
    26 = 5
    27 = 3
    28 = 7
    29 = 8
    30 = 5
    31 = 5
    32 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
